# MassWateR Censored Template: require that all parameters present in
# frecomdat also appear in censdat (qcMWRcom #70).
#
# The "Censored" sheet lists one row per parameter with a count of
# missed/censored records. This adds five new parameter rows (Water Temp,
# Sp Conductance, TP, Ammonia, E.coli) to the existing three (pH, DO,
# Nitrate), all with a default count of 0.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Censored")
$ws.Activate()

# Existing layout (row 1 = header):
#   2 pH            12
#   3 DO             1
#   4 Nitrate        0
#
# Target layout:
#   2 Water Temp     0   <- new
#   3 pH             12
#   4 DO              1
#   5 Sp Conductance  0   <- new
#   6 TP              0   <- new
#   7 Nitrate         0
#   8 Ammonia         0   <- new
#   9 E.coli          0   <- new

# Insert a blank row above pH for "Water Temp".
$ws.Rows.Item(2).Insert()

# Insert two blank rows above Nitrate (currently row 5, after the shift
# above) for "Sp Conductance" and "TP".
$ws.Rows.Item(5).Insert()
$ws.Rows.Item(5).Insert()

# Append two blank rows after Nitrate (now row 7) for "Ammonia" and
# "E.coli".
$ws.Rows.Item(8).Insert()
$ws.Rows.Item(9).Insert()

# Fill in the new rows.
$ws.Range("A2").Value2 = "Water Temp"
$ws.Range("B2").Value2 = 0

$ws.Range("A5").Value2 = "Sp Conductance"
$ws.Range("B5").Value2 = 0

$ws.Range("A6").Value2 = "TP"
$ws.Range("B6").Value2 = 0

$ws.Range("A8").Value2 = "Ammonia"
$ws.Range("B8").Value2 = 0

$ws.Range("A9").Value2 = "E.coli"
$ws.Range("B9").Value2 = 0

# Match the author's saved selection (bottom-right pane on E11).
$ws.Range("E11").Select()
